$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: update title (D4) and link (E4)
$ws.Range("D4").Value = "PyTorch 로컬이미지 로드(ImageFolder) 후, DataLoader 생성 및 CNN이미지 분류 모델 생성, 학습, 검증 성능 측정하기"
$ws.Range("E4").Value = "https://teddylee777.github.io/pytorch/pytorch-cnn-rps"

# Row 26: update title (D26)
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "http 클라이언트 프로그램 httpie와 postman 소개"
$ws.Range("E51").Value = "https://bskyvision.com/entry/http-%ED%81%B4%EB%9D%BC%EC%9D%B4%EC%96%B8%ED%8A%B8-%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%A8-httpie%EC%99%80-postman-%EC%86%8C%EA%B0%9C"
